$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F, H, I (numeric-looking) and J, K (date-looking) need a leading
# apostrophe so Excel keeps them as literal text (e.g. "760000.0", "2025-01-20")
# instead of silently converting them to a number or a date serial.
# Plain text columns (A-E, G, L, M) do not need this treatment.

# Row 5: becomes A0001 / ricardo@gmail.com / ADM_01 / transfer / PTKSKW, no nama_paket
$ws.Range("A5").Value = "A0001"
$ws.Range("B5").Value = "2025-01-20 22:03:57.412973+07:00"
$ws.Range("C5").Value = "ricardo@gmail.com"
$ws.Range("D5").Value = "ADM_01"
$ws.Range("E5").Value = "COMPLETED"
$ws.Range("F5").Value = "'760000.0"
$ws.Range("G5").Value = "transfer"
$ws.Range("H5").Value = "'760000.0"
$ws.Range("I5").Value = "'0.0"
$ws.Range("J5").Value = "'2025-01-20"
$ws.Range("K5").Value = "'2025-01-20"
$ws.Range("L5").Value = "PTKSKW"
$ws.Range("M5").ClearContents()

# Row 6: becomes A0002 / user1@gmail.com / KSR_01 / transfer / PTKSKW, no nama_paket
$ws.Range("A6").Value = "A0002"
$ws.Range("B6").Value = "2025-01-20 22:13:39.509377+07:00"
$ws.Range("C6").Value = "user1@gmail.com"
$ws.Range("D6").Value = "KSR_01"
$ws.Range("E6").Value = "COMPLETED"
$ws.Range("F6").Value = "'2280000.0"
$ws.Range("G6").Value = "transfer"
$ws.Range("H6").Value = "'2280000.0"
$ws.Range("I6").Value = "'0.0"
$ws.Range("J6").Value = "'2025-01-20"
$ws.Range("K6").Value = "'2025-01-22"
$ws.Range("L6").Value = "PTKSKW"
$ws.Range("M6").ClearContents()

# Row 7 (new): A0003 / user1@gmail.com / ADM_01 / cash / PTKSKW / Paket Wisata Ketapang
$ws.Range("A7").Value = "A0003"
$ws.Range("B7").Value = "2025-01-20 22:15:12.821899+07:00"
$ws.Range("C7").Value = "user1@gmail.com"
$ws.Range("D7").Value = "ADM_01"
$ws.Range("E7").Value = "COMPLETED"
$ws.Range("F7").Value = "'1500000.0"
$ws.Range("G7").Value = "cash"
$ws.Range("H7").Value = "'1600000.0"
$ws.Range("I7").Value = "'100000.0"
$ws.Range("J7").Value = "'2025-02-02"
$ws.Range("K7").Value = "'2025-02-04"
$ws.Range("L7").Value = "PTKSKW"
$ws.Range("M7").Value = "Paket Wisata Ketapang"
